# Auto-generated edit script: updates currentAveragePrice / LevePrice / LeveProfit
# columns across several Leve rows on multiple sheets (scheduled price-refresh run).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 34067.332
$ws.Range("J93").Value = 34067.332
$ws.Range("L93").Value = 34067.332
$ws.Range("N93").Value = -39059.332

$ws.Range("H131").Value = 4831.231
$ws.Range("I131").Value = 1311.7778
$ws.Range("J131").Value = 12750
$ws.Range("K131").Value = 3935.3334
$ws.Range("L131").Value = 38250
$ws.Range("M131").Value = 1104.6666
$ws.Range("N131").Value = -48330

$ws.Range("H132").Value = 2946.721
$ws.Range("I132").Value = 2750.0278
$ws.Range("K132").Value = 8250.0834
$ws.Range("M132").Value = -5720.0834

$ws.Range("H141").Value = 1970.4546
$ws.Range("I141").Value = 1667.5
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 5002.5
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = 177.5
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9560.048000000001
$ws.Range("I32").Value = 8432.418
$ws.Range("J32").Value = 17312.5
$ws.Range("K32").Value = 8432.418
$ws.Range("L32").Value = 17312.5
$ws.Range("M32").Value = -8145.418
$ws.Range("N32").Value = -17886.5

$ws.Range("H102").Value = 2509.7273
$ws.Range("I102").Value = 2534.3333
$ws.Range("K102").Value = 2534.3333
$ws.Range("M102").Value = -912.3332999999998

$ws.Range("H127").Value = 54990
$ws.Range("J127").Value = 54990
$ws.Range("L127").Value = 54990
$ws.Range("N127").Value = -64910

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 7193.0347
$ws.Range("I80").Value = 20209.5
$ws.Range("J80").Value = 342.26315
$ws.Range("K80").Value = 20209.5
$ws.Range("L80").Value = 342.26315
$ws.Range("M80").Value = -19211.5
$ws.Range("N80").Value = -2338.26315

$ws.Range("H83").Value = 7193.0347
$ws.Range("I83").Value = 20209.5
$ws.Range("J83").Value = 342.26315
$ws.Range("K83").Value = 101047.5
$ws.Range("L83").Value = 1711.31575
$ws.Range("M83").Value = -96055.5
$ws.Range("N83").Value = -11695.31575

$ws.Range("H134").Value = 2955.625
$ws.Range("I134").Value = 1686.9688
$ws.Range("J134").Value = 8030.25
$ws.Range("K134").Value = 5060.9064
$ws.Range("L134").Value = 24090.75
$ws.Range("M134").Value = -2525.9064
$ws.Range("N134").Value = -29160.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1100
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 1400
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1400
$ws.Range("M16").Value = -713
$ws.Range("N16").Value = -1974

$ws.Range("H107").Value = 415.16666
$ws.Range("I107").Value = 444.2
$ws.Range("J107").Value = 270
$ws.Range("K107").Value = 444.2
$ws.Range("L107").Value = 270
$ws.Range("M107").Value = 1475.8
$ws.Range("N107").Value = -4110

$ws.Range("H113").Value = 1100
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 1400
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -5740

$ws.Range("H132").Value = 14707453
$ws.Range("I132").Value = 19231828
$ws.Range("K132").Value = 57695484
$ws.Range("M132").Value = -57692954

$ws.Range("H135").Value = 54950
$ws.Range("J135").Value = 54950
$ws.Range("L135").Value = 54950
$ws.Range("N135").Value = -65090

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1314.3226
$ws.Range("I5").Value = 286.6842
$ws.Range("J5").Value = 2941.4167
$ws.Range("K5").Value = 860.0526
$ws.Range("L5").Value = 8824.250100000001
$ws.Range("M5").Value = -748.0526
$ws.Range("N5").Value = -9048.250100000001

$ws.Range("H69").Value = 1246.65
$ws.Range("J69").Value = 2464.5
$ws.Range("L69").Value = 7393.5
$ws.Range("N69").Value = -9015.5

$ws.Range("H72").Value = 1246.65
$ws.Range("J72").Value = 2464.5
$ws.Range("L72").Value = 22180.5
$ws.Range("N72").Value = -30292.5

$ws.Range("H113").Value = 733.3570999999999
$ws.Range("I113").Value = 442.93103
$ws.Range("J113").Value = 1381.2307
$ws.Range("K113").Value = 1328.79309
$ws.Range("L113").Value = 4143.6921
$ws.Range("M113").Value = 841.2069099999999
$ws.Range("N113").Value = -8483.6921

$ws.Range("H124").Value = 6600
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 6600
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 19800
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -29620

$ws.Range("H135").Value = 1314.3226
$ws.Range("I135").Value = 286.6842
$ws.Range("J135").Value = 2941.4167
$ws.Range("K135").Value = 2580.1578
$ws.Range("L135").Value = 26472.7503
$ws.Range("M135").Value = -45.15779999999995
$ws.Range("N135").Value = -31542.7503

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3935.7693
$ws.Range("I102").Value = 4136.7827
$ws.Range("J102").Value = 2394.6667
$ws.Range("K102").Value = 4136.7827
$ws.Range("L102").Value = 2394.6667
$ws.Range("M102").Value = -2514.7827
$ws.Range("N102").Value = -5638.6667

$ws.Range("H107").Value = 3903.1428
$ws.Range("J107").Value = 1299.3334
$ws.Range("L107").Value = 1299.3334
$ws.Range("N107").Value = -5139.3334

$ws.Range("H126").Value = 3726.6924
$ws.Range("I126").Value = 2494.6155
$ws.Range("J126").Value = 4958.769
$ws.Range("K126").Value = 7483.8465
$ws.Range("L126").Value = 14876.307
$ws.Range("M126").Value = -5013.8465
$ws.Range("N126").Value = -19816.307

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4336.4326
$ws.Range("I7").Value = 3952.9443
$ws.Range("J7").Value = 4699.737
$ws.Range("K7").Value = 3952.9443
$ws.Range("L7").Value = 4699.737
$ws.Range("M7").Value = -3840.9443
$ws.Range("N7").Value = -4923.737

$ws.Range("H54").Value = 39300
$ws.Range("J54").Value = 39300
$ws.Range("L54").Value = 39300
$ws.Range("N54").Value = -40588

$ws.Range("H68").Value = 1853.8667
$ws.Range("I68").Value = 1855.2727
$ws.Range("J68").Value = 1850
$ws.Range("K68").Value = 1855.2727
$ws.Range("L68").Value = 1850
$ws.Range("M68").Value = -1106.2727
$ws.Range("N68").Value = -3348

$ws.Range("H71").Value = 1853.8667
$ws.Range("I71").Value = 1855.2727
$ws.Range("J71").Value = 1850
$ws.Range("K71").Value = 9276.363499999999
$ws.Range("L71").Value = 9250
$ws.Range("M71").Value = -5532.363499999999
$ws.Range("N71").Value = -16738

$ws.Range("H122").Value = 6049.1113
$ws.Range("I122").Value = 6089.125
$ws.Range("J122").Value = 5990.909
$ws.Range("K122").Value = 18267.375
$ws.Range("L122").Value = 17972.727
$ws.Range("M122").Value = -15817.375
$ws.Range("N122").Value = -22872.727

$ws.Range("H126").Value = 4336.4326
$ws.Range("I126").Value = 3952.9443
$ws.Range("J126").Value = 4699.737
$ws.Range("K126").Value = 11858.8329
$ws.Range("L126").Value = 14099.211
$ws.Range("M126").Value = -9388.832900000001
$ws.Range("N126").Value = -19039.211

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1865.5
$ws.Range("I122").Value = 1850.3704
$ws.Range("J122").Value = 2001.6666
$ws.Range("K122").Value = 5551.1112
$ws.Range("L122").Value = 6004.9998
$ws.Range("M122").Value = -3101.1112
$ws.Range("N122").Value = -10904.9998
